function Replace-InParagraph {
    param($ParaIndex, $OldText, $NewText)
    $p = $d.Paragraphs.Item($ParaIndex).Range
    $full = $p.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Replace-InParagraph: text not found"
    }
    $absStart = $p.Start + $idx
    $absEnd = $absStart + $OldText.Length
    $sub = $d.Range($absStart, $absEnd)
    $insPt = $sub.Duplicate
    $insPt.Collapse(0)
    $insPt.InsertBefore($NewText)
    $old = $d.Range($absStart, $absEnd)
    $old.Delete()
}

$d = $word.ActiveDocument

Replace-InParagraph 3 "English" "Inglês"
Replace-InParagraph 14 "We’ll miss you at the " "Vamos sentir a sua falta "
Replace-InParagraph 14 "[EVENT NAME]" "[NOME DO EVENTO]"
Replace-InParagraph 16 "Dear " "Olá "
Replace-InParagraph 16 "[PARTNER NAME]" "[NOME DO PARCEIRO]"
Replace-InParagraph 18 "Thank you for taking the time to respond to our invitation to the upcoming " "Obrigado por ter respondido ao nosso convite para a próxima "
Replace-InParagraph 18 "[EVENT NAME]" "[NOME DO EVENTO]"
Replace-InParagraph 18 ". We were really looking forward to seeing you there." ". Gostaríamos imenso de ter a sua presença."
Replace-InParagraph 19 "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "Embora estejamos tristes por não nos podermos reunir consigo, compreendemos que por vezes surgem conflitos de agenda e outros compromissos. "
Replace-InParagraph 20 "If you’re comfortable sharing it with us, we’d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future." "Caso se sinta à vontade para o partilhar connosco, gostaríamos de saber o motivo pelo qual respondeu `"não`". Por favor, responda a este e-mail, uma vez que a sua opinião poderá ajudar-nos a melhorar os nossos processos de organização de eventos e a servi-lo melhor no futuro."
Replace-InParagraph 21 "We hope to see you at our future events. " "Esperamos contar com a sua presença nos nossos próximos eventos. "
Replace-InParagraph 22 "If you have any questions, please contact us via " "Para mais informações, contacte-nos através do "
Replace-InParagraph 22 " or " " ou "
Replace-InParagraph 23 "If you have any questions, please contact your country manager, " "Para mais questões, pode também contactar o seus gestor de parcerias "
Replace-InParagraph 23 ", at " ", em "
Replace-InParagraph 23 " or " " ou "
